# Update gh-pages generated output (丽水 event listings) with freshly
# scraped data:
#   - bump the "想去人数" counter on the first ("丽水·动漫游戏展") row
#   - insert a brand-new event ("丽水·ACG动漫游戏博览会") as the new row 3
#   - push the two events that used to sit in rows 3-4 down into rows 4-5,
#     bumping their "想去人数" counts along the way
#
# The same edit applies identically to both the "展览" sheet and the
# "全部类型" sheet (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # --- Make room: grow the table from 4 to 5 rows -----------------------
    # Row 5 doesn't exist yet. Duplicate row 4 down into row 5 first (via
    # Copy, so formatting/cell-type - e.g. the bold/bordered/centered "A"
    # column style, and the plain-text typing of the "B" date column - comes
    # along for the ride instead of Excel re-guessing types from scratch).
    $ws.Range("A4:I4").Copy($ws.Range("A5:I5"))

    # Row 3's current content (丽水·CCAC动漫游戏嘉年华) is about to move down
    # into row 4; copy its still-text "B" (start-date) cell down first so the
    # literal "2024-07-20" string isn't re-parsed into a date value.
    $ws.Range("B3").Copy($ws.Range("B4"))

    # --- Row 2: only the "想去人数" count changed --------------------------
    $ws.Range("F2").Value = 452

    # --- Row 3: brand-new event -------------------------------------------
    $ws.Range("C3").Value = "丽水·ACG动漫游戏博览会"
    $ws.Range("D3").Value = "南秦路1号望瓯·陶溪川直走200米左手边(7号楼) 望瓯陶溪川活动中心"
    $ws.Range("E3").Value = "2024.07.20 10:00-07.21 18:00"
    $ws.Range("F3").Value = 1
    $ws.Range("G3").Value = 55
    $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=86671"
    $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202405/sg6nrCrJ1717142810026.png"

    # --- Row 4: what used to be row 3 (丽水·CCAC动漫游戏嘉年华), with an
    # updated "想去人数" count. (B4 already copied from B3 above.) ---------
    $ws.Range("A4").Value = 3
    $ws.Range("C4").Value = "丽水·CCAC动漫游戏嘉年华"
    $ws.Range("D4").Value = "南环西路109号 九城宴会中心"
    $ws.Range("E4").Value = "2024.07.20 09:00-07.20 16:00"
    $ws.Range("F4").Value = 17
    $ws.Range("G4").Value = 29.9
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86306"
    $ws.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202405/4TrBjBlV1716551375116.png"

    # --- Row 5: what used to be row 4 (丽水·CCAC动漫七夕（回馈展）), with an
    # updated "想去人数" count. (B5 already copied from old B4 above.) -----
    $ws.Range("A5").Value = 4
    $ws.Range("C5").Value = "丽水·CCAC动漫七夕（回馈展）"
    $ws.Range("D5").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E5").Value = "2024.08.10 09:00-08.10 17:00"
    $ws.Range("F5").Value = 5
    $ws.Range("G5").Value = 29.9
    $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=86567"
    $ws.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"
}

Write-Output "ok"
